# Apply updated crypto price/volume data (coinranking.com scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.471.30"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "'3.385.57"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'572.68"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'137.26"
$ws.Range("E6").Value = "  +8.96%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'3.386.63"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").Value = "'0.476"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("E11").Value = "  +5.26%  "
$ws.Range("D12").Value = "'0.392"
$ws.Range("E12").Value = "  +4.82%  "
$ws.Range("D13").Value = "'3.940.40"
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "'0.0000173"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").Value = "'3.369.57"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "'25.33"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").Value = "'61.337.82"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "'14.09"
$ws.Range("E19").Value = "  +7.67%  "
$ws.Range("D20").Value = "'5.82"
$ws.Range("E20").Value = "  +4.59%  "
$ws.Range("D21").Value = "'9.44"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").Value = "'377.20"
$ws.Range("E22").Value = "  +5.87%  "
$ws.Range("D23").Value = "'0.575"
$ws.Range("D24").Value = "'3.516.11"
$ws.Range("E24").Value = "  +2.49%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'71.08"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D27").Value = "'0.0000118"
$ws.Range("E27").Value = "  +11.73%  "
$ws.Range("D28").Value = "'1.65"
$ws.Range("E28").Value = "  +12.98%  "
$ws.Range("D29").Value = "'7.78"
$ws.Range("E29").Value = "  +9.41%  "
$ws.Range("D30").Value = "'0.996"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "'8.15"
$ws.Range("E31").Value = "  +3.76%  "
$ws.Range("D32").Value = "'0.156"
$ws.Range("E32").Value = "  +5.61%  "
$ws.Range("D33").Value = "'2.15"
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'3.417.72"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").Value = "'23.57"
$ws.Range("E36").Value = "  +5.72%  "
$ws.Range("D37").Value = "'5.59"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "'7.00"
$ws.Range("E38").Value = "  +5.67%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +5.00%  "
$ws.Range("D40").Value = "'164.21"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "'0.0797"
$ws.Range("E41").Value = "  +5.70%  "
$ws.Range("D42").Value = "'0.996"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").Value = "'1.21"
$ws.Range("E43").Value = "  +9.36%  "
$ws.Range("D44").Value = "'4.42"
$ws.Range("E44").Value = "  +5.96%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'41.55"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "'0.762"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "'1.63"
$ws.Range("E47").Value = "  +6.96%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'6.99"
$ws.Range("E48").Value = "  +6.62%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'23.00"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").Value = "'23.24"
$ws.Range("E50").Value = "  +14.21%  "
$ws.Range("D51").Value = "'2.43"
$ws.Range("E51").Value = "  +14.95%  "
